# Add a new "cycle 12" column (N) to both sheets, duplicating the
# existing "cycle 11" column (M).
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(1, 14).Value = "cycle 12"
for ($r = 2; $r -le 257; $r++) {
    $ws1.Cells.Item($r, 14).Value = $ws1.Cells.Item($r, 13).Text
}

$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(1, 14).Value = "cycle 12"
for ($r = 2; $r -le 5; $r++) {
    $ws2.Cells.Item($r, 14).Value = $ws2.Cells.Item($r, 13).Text
}
